$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New section title (row 11) ---
$ws.Cells.Item(11, 1).Value = "With clamping FC layer:"
$ws.Cells.Item(11, 1).Font.Bold = $true

# --- Header row (row 13) - same labels as row 4 ---
$ws.Cells.Item(13, 2).Value = "fp32"
$ws.Cells.Item(13, 3).Value = "fp32"
$ws.Cells.Item(13, 4).Value = "ai84 quant"
$ws.Cells.Item(13, 5).Value = "ai85 quant"

# --- Sub-header row (row 14) - same labels as row 5 ---
$ws.Cells.Item(14, 1).Value = "Dataset"
$ws.Cells.Item(14, 1).Font.Italic = $true
$ws.Cells.Item(14, 2).Value = "best (verif)"
$ws.Cells.Item(14, 3).Value = "test final"
$ws.Cells.Item(14, 4).Value = "test final"
$ws.Cells.Item(14, 5).Value = "test final"

# --- Data rows 15-18 ---
$ws.Cells.Item(15, 1).Value = "MNIST"
$ws.Cells.Item(15, 2).Value = 99.4
$ws.Cells.Item(15, 3).Value = 99.6
$ws.Cells.Item(15, 4).Value = 99.5

$ws.Cells.Item(16, 1).Value = "FashionMNIST"
$ws.Cells.Item(16, 2).Value = 92.3
$ws.Cells.Item(16, 3).Value = 92.1
$ws.Cells.Item(16, 4).Value = 91.7

$ws.Cells.Item(17, 1).Value = "CIFAR-10"
$ws.Cells.Item(17, 2).Value = 82.6
$ws.Cells.Item(17, 3).Value = 82
$ws.Cells.Item(17, 4).Value = 82.4

$ws.Cells.Item(18, 1).Value = "CIFAR-10 w/bias"
$ws.Cells.Item(18, 2).Value = 82.7
$ws.Cells.Item(18, 3).Value = 82.1
$ws.Cells.Item(18, 4).Value = 31.1
$ws.Cells.Item(18, 5).Value = 81.6

# Apply the same "0.0" number format used by the existing performance
# tables to the newly added numeric cells.
$ws.Range("B15:D18").NumberFormat = "0.0"
$ws.Range("E18").NumberFormat = "0.0"

# Move the active selection to E12, like the source workbook.
$ws.Range("E12").Select() | Out-Null
